$d = $word.ActiveDocument

$replacements = @(
    @{old="15×15=225"; new="90×61=5490"},
    @{old="51×72=3672"; new="95×71=6745"},
    @{old="70×48=3360"; new="41×74=3034"},
    @{old="77×16=1232"; new="96×82=7872"},
    @{old="46×11=506"; new="50×64=3200"},
    @{old="50×35=1750"; new="15×74=1110"},
    @{old="87×74=6438"; new="97×98=9506"},
    @{old="84×32=2688"; new="74×64=4736"},
    @{old="79×80=6320"; new="93×93=8649"},
    @{old="74×18=1332"; new="99×81=8019"},
    @{old="96×49=4704"; new="60×40=2400"},
    @{old="96×39=3744"; new="16×98=1568"},
    @{old="20×48=960"; new="85×39=3315"},
    @{old="12×75=900"; new="14×53=742"},
    @{old="25×70=1750"; new="72×50=3600"},
    @{old="24×24=576"; new="95×99=9405"},
    @{old="51×15=765"; new="25×89=2225"},
    @{old="26×28=728"; new="63×21=1323"},
    @{old="26×65=1690"; new="91×40=3640"},
    @{old="77×40=3080"; new="55×88=4840"},
    @{old="19×95=1805"; new="66×71=4686"},
    @{old="63×12=756"; new="46×29=1334"},
    @{old="75×51=3825"; new="20×80=1600"},
    @{old="83×83=6889"; new="82×94=7708"},
    @{old="38×40=1520"; new="73×25=1825"}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
